$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- Build the two new border styles on sheet1's C1/D1 first ---
# C1 -> top+bottom border (reuses existing borderId=4)
$c1s1 = $ws1.Range("C1")
$c1s1.Style = "Normal"
$c1s1.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$c1s1.Borders.Item(9).LineStyle = 1    # xlEdgeBottom

# D1 -> top+right+bottom border (reuses existing borderId=5)
$d1s1 = $ws1.Range("D1")
$d1s1.Style = "Normal"
$d1s1.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$d1s1.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$d1s1.Borders.Item(9).LineStyle = 1    # xlEdgeBottom

# --- Propagate those exact formats (format-only copy, no intermediate states) ---
$ws1.Range("C1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)

$ws1.Range("D1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Text changes: anonymize "fedcore" -> "approach" in header rows ---
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Remove the stray empty inline-string cell G5 on sheet2 ---
$ws2.Range("G5").ClearContents()
